# Edit: insert a new weekly data row at row 68 of Sheet1, shifting all
# existing data rows (68-157) down by one (to 69-158), and populate the
# newly inserted row 68 with the new week's record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 68; this shifts rows 68:157 down to 69:158
# and copies formatting (e.g. the date number format on column D) from
# the row above, same as Excel's normal "Insert Sheet Rows" behavior.
$ws.Rows("68:68").Insert()

# Populate the newly inserted row 68 with the new record.
$ws.Cells.Item(68, 1).Value = 10
$ws.Cells.Item(68, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(68, 3).Value = "La Araucanía"
$ws.Cells.Item(68, 4).Value = 44482
$ws.Cells.Item(68, 5).Value = 9
$ws.Cells.Item(68, 6).Value = 100112039
$ws.Cells.Item(68, 7).Value = "Ciboulette"
$ws.Cells.Item(68, 8).Value = "Sin especificar"
$ws.Cells.Item(68, 9).Value = "Primera"
$ws.Cells.Item(68, 10).Value = 20
$ws.Cells.Item(68, 11).Value = 7000
$ws.Cells.Item(68, 12).Value = 7000
$ws.Cells.Item(68, 13).Value = 7000
$ws.Cells.Item(68, 14).Value = "$/docena de atados"
$ws.Cells.Item(68, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(68, 16).Value = 2333
$ws.Cells.Item(68, 17).Value = 3
$ws.Cells.Item(68, 18).Value = "Hortaliza"
